$d = $word.ActiveDocument
$d.Content.Find.Execute("Mecideköy", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mecidiyeköy", 2)
